# Actualización automática hashcode jue dic  6 01:55:04 CET 2018
# Updates the hash-code column (B) for specific rows in the hashcode
# worksheet, replacing stale MD5-like hash strings with updated ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value  = "8d34570a89896d1e7487e15264d4430d"
$ws.Range("B15").Value  = "0bddf56da5f82c82d65d055dd8069f4d"
$ws.Range("B29").Value  = "c9ab33bf5bace551342158f57f5fb0c5"
$ws.Range("B121").Value = "27ce3918723a74c22be7d3b4776af7d0"
$ws.Range("B126").Value = "51bbf56d85cc17f3c8cb856bf4fd262d"
$ws.Range("B133").Value = "9dbe7ba8439d5e40bc74fcda0b6edeff"
$ws.Range("B175").Value = "5da91005c11a6f40ada11d35431e6104"
$ws.Range("B191").Value = "25405b62f8f89eccdca32bc6c57b9cc6"
$ws.Range("B198").Value = "29201b2c540a545c238cd7110a9ed115"
$ws.Range("B419").Value = "2ee5add6736bc97726d8045230c25adb"
$ws.Range("B480").Value = "8a052fa960a6a06dd3c307dcbccd6d8b"
$ws.Range("B515").Value = "b6d31e86d0c877e6d1e219b37085ab4e"
$ws.Range("B547").Value = "a5555c1b1125d3fd2fcb157f929b1e11"
$ws.Range("B563").Value = "e36dde274970a017fcdcb0f19f6bba4c"
$ws.Range("B572").Value = "0751fcd52a01e68b0dea88477cc78546"
$ws.Range("B629").Value = "00d68d50c3de3d47c92bdab22d9dc903"
$ws.Range("B649").Value = "759613b2f4e599e5bbf90a4d43e40cc9"
$ws.Range("B655").Value = "6a5e3c6b8da31df5f747f3f32e2ebcf8"
$ws.Range("B733").Value = "bee7041dbfb49eb50a51ba51e5e8bca2"
$ws.Range("B862").Value = "56ad9242b497ae392e8130d0697a5abd"
